$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" '26.359.52'
$ws.Range("E2").Value = '  +0.25%  '
Set-TextCell "D3" '1.696.78'
$ws.Range("E3").Value = '  +1.10%  '
Set-TextCell "D4" '1.009'
$ws.Range("E4").Value = '  +0.02%  '
Set-TextCell "D5" '218.50'
$ws.Range("E5").Value = '  +0.02%  '
Set-TextCell "D6" '0.5424'
$ws.Range("E6").Value = '  +2.99%  '
Set-TextCell "D7" '1.009'
$ws.Range("E7").Value = '  +0.06%  '
Set-TextCell "D8" '0.2747'
$ws.Range("E8").Value = '  +1.56%  '
Set-TextCell "D9" '0.06447'
$ws.Range("E9").Value = '  -0.18%  '
Set-TextCell "D10" '21.70'
$ws.Range("E10").Value = '  -1.25%  '
Set-TextCell "D11" '0.07667'
$ws.Range("E11").Value = '  +2.03%  '
Set-TextCell "D12" '1.703.58'
$ws.Range("E12").Value = '  +2.01%  '
Set-TextCell "D13" '4.561'
$ws.Range("E13").Value = '  +0.66%  '
Set-TextCell "D14" '0.5828'
$ws.Range("E14").Value = '  +0.28%  '
Set-TextCell "D15" '0.000008401'
$ws.Range("E15").Value = '  -0.99%  '
Set-TextCell "D16" '66.50'
$ws.Range("E16").Value = '  +3.02%  '
Set-TextCell "D17" '26.416.95'
$ws.Range("E17").Value = '  +0.28%  '
Set-TextCell "D18" '4.941'
$ws.Range("E18").Value = '  +0.31%  '
Set-TextCell "D19" '1.008'
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("E20").Value = '  +0.57%  '
Set-TextCell "D21" '191.82'
$ws.Range("E21").Value = '  +1.06%  '
Set-TextCell "D22" '6.293'
$ws.Range("E22").Value = '  +1.63%  '
Set-TextCell "D23" '1.009'
$ws.Range("E23").Value = '  +0.02%  '
Set-TextCell "D24" '149.36'
$ws.Range("E24").Value = '  +2.83%  '
Set-TextCell "D25" '0.1296'
$ws.Range("E25").Value = '  +4.39%  '
Set-TextCell "D26" '7.901'
$ws.Range("E26").Value = '  +1.24%  '
Set-TextCell "D27" '15.96'
$ws.Range("E27").Value = '  +0.99%  '
$ws.Range("B28").Value = 'Hedera'
$ws.Range("C28").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell "D28" '0.06400'
$ws.Range("E28").Value = '  -1.73%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell "D29" '1.387'
$ws.Range("E29").Value = '  +2.47%  '
Set-TextCell "D30" '1.329'
$ws.Range("E30").Value = '  -0.09%  '
Set-TextCell "D31" '3.621'
$ws.Range("E31").Value = '  +0.97%  '
Set-TextCell "D32" '3.592'
$ws.Range("E32").Value = '  -0.19%  '
Set-TextCell "D33" '1.690'
$ws.Range("E33").Value = '  +1.91%  '
Set-TextCell "D34" '1.040'
$ws.Range("E34").Value = '  +0.73%  '
Set-TextCell "D35" '0.6202'
$ws.Range("E35").Value = '  -0.39%  '
Set-TextCell "D36" '2.413'
$ws.Range("E36").Value = '  +0.34%  '
Set-TextCell "D37" '2.758'
$ws.Range("E37").Value = '  +0.47%  '
Set-TextCell "D38" '0.01658'
$ws.Range("E38").Value = '  +2.31%  '
Set-TextCell "D39" '1.117.61'
$ws.Range("E39").Value = '  +0.02%  '
Set-TextCell "D40" '6.123'
$ws.Range("E40").Value = '  -5.42%  '
Set-TextCell "D41" '0.8854'
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("E42").Value = '  -0.02%  '
Set-TextCell "D43" '101.29'
$ws.Range("E43").Value = '  +0.66%  '
Set-TextCell "D44" '1.846.59'
$ws.Range("E44").Value = '  +0.78%  '
Set-TextCell "D45" '0.00000000112'
$ws.Range("E45").Value = '  +0.97%  '
Set-TextCell "D46" '57.93'
$ws.Range("E46").Value = '  +1.71%  '
Set-TextCell "D47" '8.217'
$ws.Range("E47").Value = '  +0.71%  '
Set-TextCell "D48" '1.011'
$ws.Range("E48").Value = '  +0.64%  '
Set-TextCell "D49" '0.05292'
$ws.Range("E49").Value = '  +0.23%  '
Set-TextCell "D50" '6.105'
$ws.Range("E50").Value = '  +0.18%  '
Set-TextCell "D51" '0.4301'
$ws.Range("E51").Value = '  -0.02%  '
